# Apply the crypto price/volume refresh described in the commit diff.
# Rows 17/18, 34/35 and 50/51 swap coin name+link while all D/E figures refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.217.70"
$ws.Range("E2").Value = "  +0.52%  "

$ws.Range("D3").Value = "3.022.07"
$ws.Range("E3").Value = "  -2.75%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").Value = "'556.69"
$ws.Range("E5").Value = "  -0.34%  "

$ws.Range("D6").Value = "'154.78"
$ws.Range("E6").Value = "  -4.61%  "

$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("D8").Value = "'0.560"
$ws.Range("E8").Value = "  -4.43%  "

$ws.Range("D9").Value = "3.030.71"
$ws.Range("E9").Value = "  -2.29%  "

$ws.Range("E10").Value = "  -2.31%  "

$ws.Range("D11").Value = "'6.41"
$ws.Range("E11").Value = "  -4.76%  "

$ws.Range("D12").Value = "'0.366"
$ws.Range("E12").Value = "  -2.81%  "

$ws.Range("D13").Value = "3.562.91"
$ws.Range("E13").Value = "  -2.11%  "

$ws.Range("E14").Value = "  -3.32%  "

$ws.Range("D15").Value = "63.267.59"
$ws.Range("E15").Value = "  +0.46%  "

$ws.Range("D16").Value = "'24.07"
$ws.Range("E16").Value = "  -1.92%  "

$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.034.14"
$ws.Range("E17").Value = "  -2.27%  "

$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").Value = "'0.0000150"
$ws.Range("E18").Value = "  -0.94%  "

$ws.Range("D19").Value = "'399.00"
$ws.Range("E19").Value = "  -0.42%  "

$ws.Range("D20").Value = "'5.09"
$ws.Range("E20").Value = "  -0.64%  "

$ws.Range("D21").Value = "'12.00"
$ws.Range("E21").Value = "  -2.60%  "

$ws.Range("D22").Value = "'6.65"
$ws.Range("E22").Value = "  -5.27%  "

$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = "  -0.03%  "

$ws.Range("D24").Value = "'65.32"
$ws.Range("E24").Value = "  -3.22%  "

$ws.Range("D25").Value = "'0.189"
$ws.Range("E25").Value = "  -4.66%  "

$ws.Range("E26").Value = "  -2.39%  "

$ws.Range("D27").Value = "0.0₃0985"
$ws.Range("E27").Value = "  -1.79%  "

$ws.Range("D28").Value = "'8.70"
$ws.Range("E28").Value = "  +1.06%  "

$ws.Range("D29").Value = "'0.995"
$ws.Range("E29").Value = "  -0.51%  "

$ws.Range("E30").Value = "  +0.00%  "

$ws.Range("E31").Value = "  -0.89%  "

$ws.Range("D32").Value = "'20.39"
$ws.Range("E32").Value = "  -2.06%  "

$ws.Range("D33").Value = "'162.23"
$ws.Range("E33").Value = "  +6.87%  "

$ws.Range("B34").Value = "NEARProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D34").Value = "'4.73"
$ws.Range("E34").Value = "  -1.65%  "

$ws.Range("B35").Value = "Fetch.AI"
$ws.Range("C35").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D35").Value = "'1.11"
$ws.Range("E35").Value = "  +1.48%  "

$ws.Range("D36").Value = "'6.03"
$ws.Range("E36").Value = "  -2.07%  "

$ws.Range("E37").Value = "  -0.20%  "

$ws.Range("D38").Value = "2.544.26"
$ws.Range("E38").Value = "  -5.56%  "

$ws.Range("D39").Value = "'1.59"
$ws.Range("E39").Value = "  -3.31%  "

$ws.Range("D40").Value = "'22.95"
$ws.Range("E40").Value = "  -1.22%  "

$ws.Range("D41").Value = "'3.95"
$ws.Range("E41").Value = "  -1.43%  "

$ws.Range("D42").Value = "'37.78"
$ws.Range("E42").Value = "  -0.95%  "

$ws.Range("D43").Value = "'0.670"
$ws.Range("E43").Value = "  -3.00%  "

$ws.Range("D44").Value = "'0.0599"
$ws.Range("E44").Value = "  -0.31%  "

$ws.Range("D45").Value = "'0.0250"
$ws.Range("E45").Value = "  -1.18%  "

$ws.Range("D46").Value = "'5.09"
$ws.Range("E46").Value = "  -1.51%  "

$ws.Range("D47").Value = "'0.999"
$ws.Range("E47").Value = "  -0.08%  "

$ws.Range("D48").Value = "'20.27"
$ws.Range("E48").Value = "  -1.81%  "

$ws.Range("D49").Value = "'270.01"
$ws.Range("E49").Value = "  -4.41%  "

$ws.Range("B50").Value = "WhiteBITCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D50").Value = "'10.50"
$ws.Range("E50").Value = "  +0.34%  "

$ws.Range("B51").Value = "Stellar"
$ws.Range("C51").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D51").Value = "'0.0943"
$ws.Range("E51").Value = "  -2.81%  "
